# Updates cryptocurrency price/volume figures per the latest scrape.
# (Also swaps the OKB / Filecoin rows at 32-33, which changed rank order.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '51.261.10'
$ws.Range('E2').Value = '  -1.42%  '
# Row 3
$ws.Range('D3').Value = '2.768.72'
$ws.Range('E3').Value = '  -0.42%  '
# Row 4
$ws.Range('E4').Value = '  +0.06%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '353.29'
$ws.Range('E5').Value = '  -0.85%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '107.43'
$ws.Range('E6').Value = '  -1.57%  '
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.547'
$ws.Range('E7').Value = '  -3.10%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.03%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.582'
$ws.Range('E9').Value = '  -1.62%  '
# Row 10
$ws.Range('E10').Value = '  -1.63%  '
# Row 11
$ws.Range('E11').Value = '  +3.23%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.09'
$ws.Range('E12').Value = '  +3.55%  '
# Row 13
$ws.Range('E13').Value = '  -2.48%  '
# Row 14
$ws.Range('E14').Value = '  -1.40%  '
# Row 15
$ws.Range('D15').Value = '3.203.86'
$ws.Range('E15').Value = '  -0.32%  '
# Row 16
$ws.Range('D16').Value = '2.769.34'
$ws.Range('E16').Value = '  -0.54%  '
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.922'
$ws.Range('E17').Value = '  -0.92%  '
# Row 18
$ws.Range('D18').Value = '51.227.53'
$ws.Range('E18').Value = '  -1.24%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.63'
$ws.Range('E19').Value = '  +2.94%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.09'
$ws.Range('E20').Value = '  -1.62%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.10'
$ws.Range('E21').Value = '  +0.23%  '
# Row 22
$ws.Range('D22').Value = '0.0₃0958'
$ws.Range('E22').Value = '  -1.85%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.62'
$ws.Range('E23').Value = '  -0.39%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '265.12'
$ws.Range('E24').Value = '  -3.45%  '
# Row 25
$ws.Range('E25').Value = '  -0.96%  '
# Row 26
$ws.Range('E26').Value = '  -0.02%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.90'
$ws.Range('E27').Value = '  -2.54%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.161'
$ws.Range('E28').Value = '  +12.92%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.16'
$ws.Range('E29').Value = '  +0.26%  '
# Row 30
$ws.Range('E30').Value = '  -0.67%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.58'
$ws.Range('E31').Value = '  +4.95%  '
# Row 32
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.11'
$ws.Range('E32').Value = '  +7.08%  '
# Row 33
$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '51.77'
$ws.Range('E33').Value = '  +0.14%  '
# Row 34
$ws.Range('E34').Value = '  -5.31%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.48'
$ws.Range('E35').Value = '  +3.29%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0826'
$ws.Range('E36').Value = '  -1.96%  '
# Row 37
$ws.Range('E37').Value = '  +0.06%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.11'
$ws.Range('E38').Value = '  +0.23%  '
# Row 39
$ws.Range('E39').Value = '  -2.84%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.95'
$ws.Range('E40').Value = '  -1.96%  '
# Row 41
$ws.Range('E41').Value = '  -0.57%  '
# Row 42
$ws.Range('E42').Value = '  -1.34%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '120.55'
$ws.Range('E43').Value = '  -0.76%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.04'
$ws.Range('E44').Value = '  +0.24%  '
# Row 45
$ws.Range('E45').Value = '  -2.10%  '
# Row 46
$ws.Range('D46').Value = '2.097.07'
$ws.Range('E46').Value = '  +1.85%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.23'
$ws.Range('E47').Value = '  -0.55%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.30'
$ws.Range('E48').Value = '  +1.38%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.903'
$ws.Range('E49').Value = '  -2.36%  '
# Row 50
$ws.Range('E50').Value = '  -5.30%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.31'
$ws.Range('E51').Value = '  +7.20%  '
